$wb = $excel.ActiveWorkbook

# Source sheet to clone: "31-Mar" (last existing day sheet).
$src = $wb.Worksheets.Item("31-Mar")

# --- Add "1-Apr" sheet as a copy of "31-Mar" -----------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $last)
$apr1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$apr1.Name = "1-Apr"
$apr1.Range("G10:G14").ClearContents()

# --- Add "3-Apr" sheet as a copy of "31-Mar" -----------------------------
$last2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $last2)
$apr3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$apr3.Name = "3-Apr"
$apr3.Range("G10:G14").ClearContents()

# --- View state -----------------------------------------------------------
# "31-Mar" is no longer the active tab; select all cells on it.
$src.Activate()
$src.Cells.Select()

# "3-Apr" opens scrolled so column D is the left-most visible column.
$apr3.Activate()
$apr3.Range("G14").Select()
$apr3.Application.ActiveWindow.ScrollColumn = 4

# "1-Apr" ends up as the active tab/sheet, cursor at G14.
$apr1.Activate()
$apr1.Range("G14").Select()
